$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions symbol-list update).
# Numeric-looking values (prices, percentages) are prefixed with a leading
# apostrophe so Excel stores/keeps them as text, matching the source data,
# instead of auto-converting them to numbers or percentages.
$ws.Range("D2").Value = '''307.68'
$ws.Range("E2").Value = '''0.83%'
$ws.Range("D3").Value = '''41.05'
$ws.Range("E3").Value = '''4.13%'
$ws.Range("D4").Value = '''5.124'
$ws.Range("E4").Value = '''1.94%'
$ws.Range("D5").Value = '''0.07623'
$ws.Range("E5").Value = '''-0.57%'
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = '''4.262'
$ws.Range("E6").Value = '''0.38%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = '''1.619'
$ws.Range("E7").Value = '''0.84%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = '''2.508'
$ws.Range("E8").Value = '''1.86%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '''0.9017'
$ws.Range("E9").Value = '''2.14%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.1091'
$ws.Range("E10").Value = '''12.55%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1776'
$ws.Range("E11").Value = '''3.25%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.09158'
$ws.Range("E12").Value = '''2.96%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.04196'
$ws.Range("E13").Value = '''-6.49%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.1051'
$ws.Range("E14").Value = '''-0.51%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001254'
$ws.Range("E15").Value = '''1.18%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.005792'
$ws.Range("E16").Value = '''-2.46%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.357'
$ws.Range("E17").Value = '''0.15%'
$ws.Range("D19").Value = '''6.556'
$ws.Range("E19").Value = '''-6.92%'
$ws.Range("E20").Value = '''0.99%'
$ws.Range("E21").Value = '''-12.16%'
$ws.Range("D22").Value = '''0.04073'
$ws.Range("E22").Value = '''-3.13%'
$ws.Range("E23").Value = '''2.21%'
$ws.Range("D24").Value = '''0.004002'
$ws.Range("E24").Value = '''-1.49%'
$ws.Range("E25").Value = '''6.44%'
$ws.Range("D38").Value = '''0.02389'
$ws.Range("E38").Value = '''2.52%'
$ws.Range("E39").Value = '''1.05%'
$ws.Range("D40").Value = '''0.007768'
$ws.Range("E40").Value = '''-1.98%'
$ws.Range("E41").Value = '''-1.74%'
$ws.Range("D42").Value = '''0.006878'
$ws.Range("E42").Value = '''6.28%'
$ws.Range("E43").Value = '''-1.86%'
$ws.Range("D44").Value = '''0.008544'
$ws.Range("E44").Value = '''-1.33%'
$ws.Range("E45").Value = '''1.46%'
$ws.Range("D46").Value = '''0.00006887'
$ws.Range("E46").Value = '''5.33%'
$ws.Range("E47").Value = '''-0.15%'
$ws.Range("D48").Value = '''0.01182'
$ws.Range("E48").Value = '''250.54%'
$ws.Range("D49").Value = '''0.004200'
$ws.Range("E49").Value = '''-40.05%'
$ws.Range("D50").Value = '''0.00002100'
$ws.Range("E50").Value = '''-0.15%'
$ws.Range("D51").Value = '''0.0002000'
$ws.Range("E51").Value = '''-0.15%'
